# Natmi following Dr Hou advice:
# Rework the Efna5-Epha7 LR-pair sheet so the "Target cluster" column (D)
# spans the full FAPs / sCs / ECs cluster set (adding the new "ECs" cluster),
# updating the corresponding computed metrics, and adding the two new rows
# (sCs -> ECs and the filled-in sCs x {FAPs,sCs} combinations) this produces.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.666083666666667
$ws.Range("H2").Value = 4.998251
$ws.Range("I2").Value = 0.6125276070882968
$ws.Range("J2").Value = 0.6125276070882968
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1532026666666667
$ws.Range("N2").Value = 0.459608
$ws.Range("O2").Value = 0.03935316038222934
$ws.Range("P2").Value = 0.03935316038222934
$ws.Range("Q2").Value = 0.2552484606231111
$ws.Range("R2").Value = 2.297236145608
$ws.Range("S2").Value = 0.0241048971602889
$ws.Range("T2").Value = 0.0241048971602889

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.666083666666667
$ws.Range("H3").Value = 4.998251
$ws.Range("I3").Value = 0.6125276070882968
$ws.Range("J3").Value = 0.6125276070882968
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1693566666666667
$ws.Range("N3").Value = 0.50807
$ws.Range("O3").Value = 0.04350263745496
$ws.Range("P3").Value = 0.04350263745496
$ws.Range("Q3").Value = 0.2821623761744445
$ws.Range("R3").Value = 2.53946138557
$ws.Range("S3").Value = 0.02664656642231637
$ws.Range("T3").Value = 0.02664656642231637

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha7"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.666083666666667
$ws.Range("H4").Value = 4.998251
$ws.Range("I4").Value = 0.6125276070882968
$ws.Range("J4").Value = 0.6125276070882968
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.570461333333333
$ws.Range("N4").Value = 10.711384
$ws.Range("O4").Value = 0.9171442021628107
$ws.Range("P4").Value = 0.9171442021628107
$ws.Range("Q4").Value = 5.948687309931554
$ws.Range("R4").Value = 53.53818578938399
$ws.Range("S4").Value = 0.5617761435056915
$ws.Range("T4").Value = 0.5617761435056915

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.053930333333333
$ws.Range("H5").Value = 3.161791
$ws.Range("I5").Value = 0.3874723929117032
$ws.Range("J5").Value = 0.3874723929117031
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1532026666666667
$ws.Range("N5").Value = 0.459608
$ws.Range("O5").Value = 0.03935316038222934
$ws.Range("P5").Value = 0.03935316038222934
$ws.Range("Q5").Value = 0.1614649375475556
$ws.Range("R5").Value = 1.453184437928
$ws.Range("S5").Value = 0.01524826322194044
$ws.Range("T5").Value = 0.01524826322194043

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.053930333333333
$ws.Range("H6").Value = 3.161791
$ws.Range("I6").Value = 0.3874723929117032
$ws.Range("J6").Value = 0.3874723929117031
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1693566666666667
$ws.Range("N6").Value = 0.50807
$ws.Range("O6").Value = 0.04350263745496
$ws.Range("P6").Value = 0.04350263745496
$ws.Range("Q6").Value = 0.1784901281522223
$ws.Range("R6").Value = 1.60641115337
$ws.Range("S6").Value = 0.01685607103264364
$ws.Range("T6").Value = 0.01685607103264363

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha7"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.053930333333333
$ws.Range("H7").Value = 3.161791
$ws.Range("I7").Value = 0.3874723929117032
$ws.Range("J7").Value = 0.3874723929117031
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.570461333333333
$ws.Range("N7").Value = 10.711384
$ws.Range("O7").Value = 0.9171442021628107
$ws.Range("P7").Value = 0.9171442021628107
$ws.Range("Q7").Value = 3.763017503193777
$ws.Range("R7").Value = 33.86715752874399
$ws.Range("S7").Value = 0.3553680586571191
$ws.Range("T7").Value = 0.355368058657119
